# Update cryptocurrency price/volume data to the latest scraped snapshot.
# (commit: "Updated cryptos list on Sun Jul 23 13:30:23 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.955.53"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").Value = "1.878.34"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'0.7406"
$ws.Range("E5").Value = "  -3.90%  "

$ws.Range("D6").Value = "'243.11"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").Value = "'0.3153"
$ws.Range("E8").Value = "  +1.11%  "

$ws.Range("D9").Value = "'0.07204"
$ws.Range("E9").Value = "  +0.53%  "

$ws.Range("D10").Value = "'24.64"
$ws.Range("E10").Value = "  -3.91%  "

$ws.Range("D11").Value = "'0.08344"
$ws.Range("E11").Value = "  -2.94%  "

$ws.Range("D12").Value = "'0.7506"
$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.405"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.863.80"
$ws.Range("E14").Value = "  -1.87%  "

$ws.Range("D15").Value = "'92.48"
$ws.Range("E15").Value = "  -1.22%  "

$ws.Range("D16").Value = "29.985.59"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("D17").Value = "'6.102"
$ws.Range("E17").Value = "  -0.85%  "

$ws.Range("D18").Value = "'248.41"
$ws.Range("E18").Value = "  +1.58%  "

$ws.Range("D19").Value = "'13.56"
$ws.Range("E19").Value = "  -1.54%  "

$ws.Range("D20").Value = "'0.000007854"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("D22").Value = "2.143.76"
$ws.Range("E22").Value = "  -2.20%  "

$ws.Range("D23").Value = "'8.019"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("E25").Value = "  -5.66%  "

$ws.Range("D26").Value = "'9.279"
$ws.Range("E26").Value = "  -1.08%  "

$ws.Range("D27").Value = "'164.93"
$ws.Range("E27").Value = "  +1.45%  "

$ws.Range("D28").Value = "'18.69"
$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").Value = "'2.035"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").Value = "'1.510"
$ws.Range("E30").Value = "  +3.40%  "

$ws.Range("D31").Value = "'4.599"
$ws.Range("E31").Value = "  +1.97%  "

$ws.Range("D32").Value = "'1.537"
$ws.Range("E32").Value = "  +0.18%  "

$ws.Range("D33").Value = "'4.269"
$ws.Range("E33").Value = "  +4.26%  "

$ws.Range("D34").Value = "'0.05321"
$ws.Range("E34").Value = "  -2.63%  "

$ws.Range("E35").Value = "  -0.40%  "

$ws.Range("D36").Value = "'0.7495"
$ws.Range("E36").Value = "  +0.69%  "

$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("D38").Value = "'2.695"
$ws.Range("E38").Value = "  -0.14%  "

$ws.Range("D39").Value = "'0.01968"
$ws.Range("E39").Value = "  +0.77%  "

$ws.Range("D40").Value = "'2.757"
$ws.Range("E40").Value = "  -0.85%  "

$ws.Range("D41").Value = "'0.4547"
$ws.Range("E41").Value = "  +1.78%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'6.142"
$ws.Range("E42").Value = "  +0.98%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.107.23"
$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("D44").Value = "'72.22"
$ws.Range("E44").Value = "  -1.12%  "

$ws.Range("D45").Value = "'0.8588"
$ws.Range("E45").Value = "  +0.82%  "

$ws.Range("D46").Value = "'104.32"
$ws.Range("E46").Value = "  +1.80%  "

$ws.Range("E47").Value = "  +0.34%  "

$ws.Range("D48").Value = "'1.857"
$ws.Range("E48").Value = "  -0.37%  "

$ws.Range("D49").Value = "'7.617"
$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("D50").Value = "'9.501"
$ws.Range("E50").Value = "  -2.80%  "

$ws.Range("D51").Value = "2.039.81"
$ws.Range("E51").Value = "  -1.35%  "
